$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# 4 Bölümün cetvel de yoksa devam etme özelliği eklendi.
# Update DateNumber_1 / DateNumber_2 from 14 to 23
$ws.Range("B8").Value = 23
$ws.Range("B9").Value = 23
$ws.Range("B11").Value = 23
$ws.Range("B12").Value = 23

# Clear DateNumber_3 / DateNumber_4 values (keep formatting)
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("B18").ClearContents()
